$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.03530768103927
    "D2" = 1.04204750479483
    "E2" = 0.992614727750844
    "F2" = 1.049452705300672
    "I2" = 1.038090244184041
    "J2" = 1.040422168627361
    "K2" = 1.044825296074835
    "L2" = 0.9955398523335997
    "M2" = 1.052209723029879
    "N2" = 1.041899687444903
    "B3" = 1.02
    "C3" = 1.036104413087579
    "D3" = 1.042664554177294
    "E3" = 0.9936372048519299
    "F3" = 1.050249332629033
    "I3" = 1.038253507618088
    "J3" = 1.040863166113219
    "K3" = 1.045253605525628
    "L3" = 0.9963617723202687
    "M3" = 1.052818639721772
    "N3" = 1.042341311197774
    "B4" = 1.02
    "C4" = 1.036620511047361
    "D4" = 1.043064341085981
    "E4" = 0.9942998659930998
    "F4" = 1.050765791241896
    "I4" = 1.038358242492089
    "J4" = 1.041148397869035
    "K4" = 1.04553057108777
    "L4" = 0.9968940712668347
    "M4" = 1.053212979461426
    "N4" = 1.0426269480154
    "B5" = 1.02
    "C5" = 1.036837610696118
    "D5" = 1.043232532877899
    "E5" = 0.994578699834602
    "F5" = 1.050983144789395
    "I5" = 1.038402055061362
    "J5" = 1.041268278533691
    "K5" = 1.045646963149351
    "L5" = 0.9971179600053012
    "M5" = 1.053378837173265
    "N5" = 1.042746998924349
    "B6" = 1.02
    "C6" = 1.036874070375423
    "D6" = 1.04326078008369
    "E6" = 0.994625531979634
    "F6" = 1.051019653088884
    "I6" = 1.038409398581711
    "J6" = 1.041288405183068
    "K6" = 1.045666503254264
    "L6" = 0.9971555583673455
    "M6" = 1.053406689895527
    "N6" = 1.042767154155878
    "B7" = 1.02
    "C7" = 1.036623411426778
    "D7" = 1.043066587998334
    "E7" = 0.994303590798249
    "F7" = 1.050768694612395
    "I7" = 1.038358828775158
    "J7" = 1.041149999842853
    "K7" = 1.045532126499787
    "L7" = 0.9968970624462089
    "M7" = 1.053215195358745
    "N7" = 1.042628552264204
    "B8" = 1.02
    "C8" = 1.035576823545311
    "D8" = 1.042255931843213
    "E8" = 0.9929600610674297
    "F8" = 1.04972172354092
    "I8" = 1.038145607294648
    "J8" = 1.040571230504723
    "K8" = 1.044970081578723
    "L8" = 0.9958175282591056
    "M8" = 1.052415439942895
    "N8" = 1.042048961007228
    "B9" = 1.02
    "C9" = 1.033736962001807
    "D9" = 1.040831475577501
    "E9" = 0.9906006454969559
    "F9" = 1.047884480355137
    "I9" = 1.037762965356659
    "J9" = 1.039550479339058
    "K9" = 1.04397837259569
    "L9" = 0.9939188001724441
    "M9" = 1.051008770660329
    "N9" = 1.041026760257826
    "B10" = 1.02
    "C10" = 1.032513424714287
    "D10" = 1.039884653317156
    "E10" = 0.989033133672735
    "F10" = 1.046664920116627
    "I10" = 1.037503262102882
    "J10" = 1.038869460296474
    "K10" = 1.043316429658141
    "L10" = 0.9926553831429383
    "M10" = 1.050072841994964
    "N10" = 1.040344774090087
    "B11" = 1.02
    "C11" = 1.031984362184684
    "D11" = 1.039475359875725
    "E11" = 0.988355674866747
    "F11" = 1.046138112666359
    "I11" = 1.037389723684894
    "J11" = 1.038574464907409
    "K11" = 1.043029627328618
    "L11" = 0.9921088820399291
    "M11" = 1.049668035004611
    "N11" = 1.040049359773734
    "B12" = 1.02
    "C12" = 1.031787957352354
    "D12" = 1.03932343515219
    "E12" = 0.9881042295826724
    "F12" = 1.04594262586412
    "I12" = 1.037347388141231
    "J12" = 1.038464875323858
    "K12" = 1.042923070991078
    "L12" = 0.9919059725120875
    "M12" = 1.049517741954692
    "N12" = 1.039939614560404
    "B13" = 1.02
    "C13" = 1.03183008173377
    "D13" = 1.039356018735447
    "E13" = 0.9881581567098651
    "F13" = 1.045984549681942
    "I13" = 1.037356476592828
    "J13" = 1.03848838332656
    "K13" = 1.042945928804925
    "L13" = 0.9919494934313052
    "M13" = 1.049549977115421
    "N13" = 1.039963155947166
    "B14" = 1.02
    "C14" = 1.031968124993658
    "D14" = 1.039462799576055
    "E14" = 0.9883348863814464
    "F14" = 1.046121949718169
    "I14" = 1.037386227522765
    "J14" = 1.038565406499947
    "K14" = 1.043020819854369
    "L14" = 0.9920921077337197
    "M14" = 1.049655610290279
    "N14" = 1.040040288502294
    "B15" = 1.02
    "C15" = 1.03205319296941
    "D15" = 1.039528604757453
    "E15" = 0.9884438009545853
    "F15" = 1.046206632044901
    "I15" = 1.037404536565153
    "J15" = 1.038612861049369
    "K15" = 1.043066959405393
    "L15" = 0.9921799884222134
    "M15" = 1.049720703741723
    "N15" = 1.040087810442619
    "B16" = 1.02
    "C16" = 1.032548552536717
    "D16" = 1.039911831407645
    "E16" = 0.9890781214508737
    "F16" = 1.046699909533547
    "I16" = 1.037510774471224
    "J16" = 1.038889035970729
    "K16" = 1.043335460178082
    "L16" = 0.9926916645766087
    "M16" = 1.050099717463655
    "N16" = 1.040364377564046
    "B17" = 1.02
    "C17" = 1.032859477138005
    "D17" = 1.040152404586392
    "E17" = 0.989476357848556
    "F17" = 1.047009671023855
    "I17" = 1.037577124672016
    "J17" = 1.039062244772601
    "K17" = 1.043503837280043
    "L17" = 0.9930127773699352
    "M17" = 1.050337586148259
    "N17" = 1.040537832342284
    "B18" = 1.02
    "C18" = 1.033040905205294
    "D18" = 1.04029279299518
    "E18" = 0.9897087662937556
    "F18" = 1.047190472072985
    "I18" = 1.037615720867714
    "J18" = 1.039163263868699
    "K18" = 1.043602031657695
    "L18" = 0.9932001317071769
    "M18" = 1.050476374837291
    "N18" = 1.040638994897084
    "B19" = 1.02
    "C19" = 1.033102779477343
    "D19" = 1.040340672976092
    "E19" = 0.9897880325774034
    "F19" = 1.047252141257484
    "I19" = 1.037628863383123
    "J19" = 1.039197706920566
    "K19" = 1.043635510454746
    "L19" = 0.9932640239640975
    "M19" = 1.050523705588628
    "N19" = 1.040673486862036
    "B20" = 1.02
    "C20" = 1.032826110513245
    "D20" = 1.040126586505852
    "E20" = 0.9894336180360679
    "F20" = 1.04697642387057
    "I20" = 1.03757001675341
    "J20" = 1.039043662201325
    "K20" = 1.043485773773279
    "L20" = 0.9929783193494215
    "M20" = 1.050312060535253
    "N20" = 1.040519223381624
    "B21" = 1.02
    "C21" = 1.031927471572982
    "D21" = 1.039431352363904
    "E21" = 0.9882828385668249
    "F21" = 1.046081483489052
    "I21" = 1.03737747109047
    "J21" = 1.038542725503129
    "K21" = 1.04299876699116
    "L21" = 0.9920501090198102
    "M21" = 1.049624501998823
    "N21" = 1.040017575295859
    "B22" = 1.02
    "C22" = 1.031363113659649
    "D22" = 1.038994839986057
    "E22" = 0.9875604150241495
    "F22" = 1.045519916291179
    "I22" = 1.037255471213699
    "J22" = 1.038227679649064
    "K22" = 1.042692421472632
    "L22" = 0.9914670000341481
    "M22" = 1.049192614194928
    "N22" = 1.039702082040544
    "B23" = 1.02
    "C23" = 1.031662228060879
    "D23" = 1.039226185106498
    "E23" = 0.9879432794643023
    "F23" = 1.045817507072712
    "I23" = 1.037320234401575
    "J23" = 1.038394699232299
    "K23" = 1.042854834321019
    "L23" = 0.991776070289318
    "M23" = 1.049421526881912
    "N23" = 1.039869338810746
    "B24" = 1.02
    "C24" = 1.032841187235547
    "D24" = 1.040138252380999
    "E24" = 0.9894529299347244
    "F24" = 1.046991446450076
    "I24" = 1.037573228839634
    "J24" = 1.039052058898536
    "K24" = 1.043493935947578
    "L24" = 0.9929938892766442
    "M24" = 1.050323594325585
    "N24" = 1.040527632003109
    "B25" = 1.02
    "C25" = 1.03421208274889
    "D25" = 1.041199243510047
    "E25" = 0.9912096547607049
    "F25" = 1.04835853211739
    "I25" = 1.037862702921495
    "J25" = 1.039814464981151
    "K25" = 1.044234900626351
    "L25" = 0.9944092447426414
    "M25" = 1.051372109586624
    "N25" = 1.041291120789811
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}